# Scheduled-runner price refresh: update cached market-board pricing
# columns (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# for the leves whose prices moved since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 106: Making Your Mark | Enchanted Palladium Ink
$ws.Range("H106").Value = 25643652
$ws.Range("I106").Value = 30304770
$ws.Range("K106").Value = 30304770
$ws.Range("M106").Value = -30304139
# Row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 25401.062
$ws.Range("I113").Value = 3844.0588
$ws.Range("K113").Value = 3844.0588
$ws.Range("M113").Value = -590.0587999999998
# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 116899.664
$ws.Range("I132").Value = 283150.12
$ws.Range("K132").Value = 849450.36
$ws.Range("M132").Value = -846920.36
# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 4813.654
$ws.Range("I137").Value = 1827.8
$ws.Range("K137").Value = 5483.4
$ws.Range("M137").Value = -2933.4
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5170.0586
$ws.Range("I138").Value = 1529.0625
$ws.Range("J138").Value = 6290.365
$ws.Range("K138").Value = 4587.1875
$ws.Range("L138").Value = 18871.095
$ws.Range("M138").Value = 552.8125
$ws.Range("N138").Value = -29151.095

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 23366.23
$ws.Range("I61").Value = 19072
$ws.Range("J61").Value = 28376.166
$ws.Range("K61").Value = 19072
$ws.Range("L61").Value = 28376.166
$ws.Range("M61").Value = -18860
$ws.Range("N61").Value = -28800.166
# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 596750.9399999999
$ws.Range("I102").Value = 1246538.9
$ws.Range("K102").Value = 1246538.9
$ws.Range("M102").Value = -1244916.9
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 23366.23
$ws.Range("I136").Value = 19072
$ws.Range("J136").Value = 28376.166
$ws.Range("K136").Value = 57216
$ws.Range("L136").Value = 85128.49800000001
$ws.Range("M136").Value = -54666
$ws.Range("N136").Value = -90228.49800000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 7924.6665
$ws.Range("J20").Value = 6987
$ws.Range("L20").Value = 6987
$ws.Range("N20").Value = -7481
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 1890.3
$ws.Range("I86").Value = 2968.6667
$ws.Range("J86").Value = 1428.1428
$ws.Range("K86").Value = 2968.6667
$ws.Range("L86").Value = 1428.1428
$ws.Range("M86").Value = -1845.6667
$ws.Range("N86").Value = -3674.1428
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 1890.3
$ws.Range("I89").Value = 2968.6667
$ws.Range("J89").Value = 1428.1428
$ws.Range("K89").Value = 14843.3335
$ws.Range("L89").Value = 7140.714
$ws.Range("M89").Value = -9227.333500000001
$ws.Range("N89").Value = -18372.714
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2551.4
$ws.Range("I134").Value = 2063.6316
$ws.Range("J134").Value = 5199.2856
$ws.Range("K134").Value = 6190.8948
$ws.Range("L134").Value = 15597.8568
$ws.Range("M134").Value = -3655.8948
$ws.Range("N134").Value = -20667.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 967.5
$ws.Range("I22").Value = 710
$ws.Range("K22").Value = 710
$ws.Range("M22").Value = -360
# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 2018.7646
$ws.Range("I58").Value = 1958.7273
$ws.Range("K58").Value = 1958.7273
$ws.Range("M58").Value = -1755.7273
# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 5463.4
$ws.Range("I122").Value = 2927.4
$ws.Range("J122").Value = 7999.4
$ws.Range("K122").Value = 8782.200000000001
$ws.Range("L122").Value = 23998.2
$ws.Range("M122").Value = -6332.200000000001
$ws.Range("N122").Value = -28898.2
# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 15881654
$ws.Range("I132").Value = 16670737
$ws.Range("K132").Value = 50012211
$ws.Range("M132").Value = -50009681
# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2424.4517
$ws.Range("I134").Value = 2424.4517
$ws.Range("K134").Value = 7273.355100000001
$ws.Range("M134").Value = -4738.355100000001
# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 2018.7646
$ws.Range("I136").Value = 1958.7273
$ws.Range("K136").Value = 5876.1819
$ws.Range("M136").Value = -3326.1819

$ws = $wb.Worksheets.Item("CUL")
# Row 4: In Hot Water | Boiled Egg
$ws.Range("H4").Value = 49597080
$ws.Range("I4").Value = 1400351.8
$ws.Range("K4").Value = 4201055.4
$ws.Range("M4").Value = -4200943.4
# Row 37: I Love Lamprey | Eel Pie
$ws.Range("H37").Value = 149000
$ws.Range("J37").Value = 149000
$ws.Range("L37").Value = 447000
$ws.Range("N37").Value = -447224
# Row 80: Saucy for a Suitor | Hollandaise Sauce
$ws.Range("H80").Value = 5843.1875
$ws.Range("J80").Value = 5843.1875
$ws.Range("L80").Value = 17529.5625
$ws.Range("N80").Value = -19401.5625
# Row 83: Saved by the Sauce (L) | Hollandaise Sauce
$ws.Range("H83").Value = 5843.1875
$ws.Range("J83").Value = 5843.1875
$ws.Range("L83").Value = 52588.6875
$ws.Range("N83").Value = -61948.6875
# Row 101: No Othard Choice | Egg Foo Young
$ws.Range("H101").Value = 10995.143
$ws.Range("J101").Value = 10995.143
$ws.Range("L101").Value = 32985.429
$ws.Range("N101").Value = -37853.429
# Row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 1059.8
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1059.8
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3179.4
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7019.4
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 8335825.5
$ws.Range("J131").Value = 5749676.5
$ws.Range("L131").Value = 17249029.5
$ws.Range("N131").Value = -17259109.5
# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 2519.889
$ws.Range("I132").Value = 1649.8334
$ws.Range("K132").Value = 14848.5006
$ws.Range("M132").Value = -12318.5006

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 2981740.5
$ws.Range("I70").Value = 3973779.5
$ws.Range("K70").Value = 3973779.5
$ws.Range("M70").Value = -3973509.5
# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 2981740.5
$ws.Range("I73").Value = 3973779.5
$ws.Range("K73").Value = 3973779.5
$ws.Range("M73").Value = -3972843.5
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2203521.2
$ws.Range("I122").Value = 2203521.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6610563.600000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6608113.600000001
$ws.Range("N122").ClearContents()
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 4326.3477
$ws.Range("I126").Value = 2385.7144
$ws.Range("J126").Value = 5175.375
$ws.Range("K126").Value = 7157.1432
$ws.Range("L126").Value = 15526.125
$ws.Range("M126").Value = -4687.1432
$ws.Range("N126").Value = -20466.125

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 3842.2239
$ws.Range("I132").Value = 2868.08
$ws.Range("K132").Value = 8604.24
$ws.Range("M132").Value = -6074.24

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 1993669.5
$ws.Range("I81").Value = 1491365.9
$ws.Range("J81").Value = 2998276.5
$ws.Range("K81").Value = 2982731.8
$ws.Range("L81").Value = 5996553
$ws.Range("M81").Value = -2981670.8
$ws.Range("N81").Value = -5998675
# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 1993669.5
$ws.Range("I84").Value = 1491365.9
$ws.Range("J84").Value = 2998276.5
$ws.Range("K84").Value = 14913659
$ws.Range("L84").Value = 29982765
$ws.Range("M84").Value = -14908355
$ws.Range("N84").Value = -29993373
